$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated power-loss results for the 380 kV line case (Case_0_20)
# Columns: C,D,E,F,G,I,L,N,O for rows 2-25
$colLetters = @("C","D","E","F","G","I","L","N","O")

$data = @{
    2 = @(0.1820850610631197, 0.0213911499970223, 0.1309150718809882, 1.21103156460741, 0.00241949403477844, 1.163511203300033, 0.1758153550338051, 2.873581687696344, 4.266805778577407)
    3 = @(0.180673226964899, 0.02150530346882995, 0.1282938749517548, 1.16187534105795, 0.002424323654527597, 1.12408729910365, 0.1706748245311402, 2.562605684679454, 4.106082441381545)
    4 = @(0.1798937303394439, 0.02158405416785669, 0.1267606846778939, 1.132458981401882, 0.002427448008595435, 1.100605263042169, 0.1676285239038577, 2.371325805375761, 4.010087063132914)
    5 = @(0.1795980346744983, 0.02161831890285981, 0.1261550105877767, 1.120663222324197, 0.002428761312759271, 1.091217549903931, 0.1664146691013144, 2.293303068607429, 3.971641334929245)
    6 = @(0.1795502602560788, 0.02162413966502186, 0.126055591635005, 1.118716094122945, 0.002428981812177232, 1.089669662163672, 0.1662147690242648, 2.280343261403573, 3.96529801887624)
    7 = @(0.1798896536076739, 0.02158450748133234, 0.1267524390294206, 1.132299124909693, 0.002427465557627493, 1.100477923543778, 0.1676120420978577, 2.370273851395496, 4.009565849089768)
    8 = @(0.1815801154143628, 0.0214287099896282, 0.1299954251434094, 1.19392299650508, 0.002421126380007168, 1.149767023358507, 0.1740199911302085, 2.766433886209654, 4.210828249189092)
    9 = @(0.1855895181492144, 0.02119212000829762, 0.1369628760972219, 1.320892216726264, 0.00240995025057349, 1.252211094437143, 0.187465320791631, 3.540180268007646, 4.627008972195881)
    10 = @(0.1889607340252581, 0.0210606377703364, 0.1424575016859322, 1.417989857857179, 0.002402495546409715, 1.331071227906264, 0.1978900075677643, 4.10623028343673, 4.946156478208536)
    11 = @(0.1905873037923556, 0.02101009120175945, 0.1450398561336996, 1.463008112900809, 0.002399266604323518, 1.367742183416638, 0.2027534419483175, 4.363110593465422, 5.094311143175617)
    12 = @(0.191216647170549, 0.02099228925552055, 0.1460297154526415, 1.480178544764243, 0.00239806707645935, 1.381744171662575, 0.2046126905016621, 4.460285735714251, 5.15084518216895)
    13 = @(0.1910805105778905, 0.02099606356556194, 0.1458159977282207, 1.476475097936799, 0.002398324386342494, 1.378723442314637, 0.2042114842439844, 4.439361943450422, 5.138650335089608)
    14 = @(0.1906388116128284, 0.02100859975131542, 0.1451210520098343, 1.464418263595377, 0.00239916745413159, 1.368891815708324, 0.202906050452853, 4.371107314139522, 5.098953572300616)
    15 = @(0.1903700036573781, 0.02101645309753053, 0.1446969395863604, 1.457049153434156, 0.00239968687528944, 1.362884727316214, 0.2021087272071327, 4.329286057409945, 5.074694405999821)
    16 = @(0.1888563071406253, 0.02106412810813652, 0.142290411126659, 1.415064963122234, 0.002402709819187761, 1.328690808459214, 0.1975746216708387, 4.089429168003562, 4.936534349568376)
    17 = @(0.1879515372129248, 0.02109575287444443, 0.1408353445562867, 1.389526968570152, 0.002404605759666012, 1.307918716518429, 0.1948242369593203, 3.94211849063862, 4.852541587904739)
    18 = @(0.1874398882641088, 0.02111481429403383, 0.1400062201814407, 1.37491794862629, 0.002405711533308187, 1.296046131035396, 0.1932536868446704, 3.857331695637754, 4.80451077513834)
    19 = @(0.1872681549195505, 0.02112141767955933, 0.1397268290449105, 1.36998525252784, 0.002406088557027975, 1.292039127776562, 0.1927238788912575, 3.828614786364199, 4.788296280678196)
    20 = @(0.1880469456630607, 0.02109229610164576, 0.1409894320031029, 1.392237269412419, 0.002404402353337351, 1.310122175544663, 0.1951158393735142, 3.957806003280837, 4.861453797260253)
    21 = @(0.1907681855161485, 0.02100488117778099, 0.1453248487924981, 1.467956301613981, 0.002398919195752315, 1.3717764612813, 0.2032890099561371, 4.391158149571254, 5.110601745461395)
    22 = @(0.192624771782107, 0.02095555810053185, 0.148228142127941, 1.518160476400283, 0.002395470818995049, 1.412744699467538, 0.2087331380031543, 4.673791817957863, 5.275948717082031)
    23 = @(0.1916267220930337, 0.02098116610948608, 0.1466721859508979, 1.491299570416487, 0.002397298955384617, 1.390817253235113, 0.2058180784953976, 4.523002190001307, 5.187468660132538)
    24 = @(0.1880037849810066, 0.02109385616825676, 0.1409197459848599, 1.391011714737942, 0.002404494264198968, 1.309125775089342, 0.1949839726598555, 3.950713976768498, 4.857423784822913)
    25 = @(0.1844303053047298, 0.02124871243840332, 0.1350123955114313, 1.285879800927859, 0.002412840234357405, 1.223870829721307, 0.1837328012529582, 3.331249627311138, 4.512091093311199)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $cellRef = "$($colLetters[$i])$row"
        $ws.Range($cellRef).Value = $vals[$i]
    }
}

Write-Output "Updated $($data.Keys.Count) rows across columns $($colLetters -join ',')"